$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'36.153.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.23%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'2.013.24"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.66%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'253.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.32%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'0.642"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.55%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'62.36"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +12.08%  "
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("E8").Value = "'  -0.03%  "
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'59.08"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -7.36%  "
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.370"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.56%  "
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.0749"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.14%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("E12").Value = "'  -1.66%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'0.918"
$ws.Range("D13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'14.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.66%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'2.301.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.81%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'5.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.28%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'19.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +12.69%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'2.014.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.50%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'36.090.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.14%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'72.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.23%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'0.0₃0859"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.28%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'5.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.54%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'233.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.71%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'2.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +20.19%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "'  +0.14%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  -2.54%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'9.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +3.30%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'164.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.18%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'19.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.94%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("E30").Value = "'  -1.00%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'5.15"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +3.00%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'1.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.69%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'0.109"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +25.07%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'0.0604"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.59%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'4.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.59%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'2.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +11.99%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("E37").Value = "'  -0.09%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("E38").Value = "'  -1.21%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'5.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +16.77%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.106"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +17.78%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'1.22"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.79%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'2.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.02%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.0216"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.32%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("E44").Value = "'  +1.29%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'16.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +3.99%  "
$ws.Range("E45").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'1.411.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.58%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'2.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +4.23%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'2.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.96%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'47.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +2.99%  "
$ws.Range("E51").Style = "Normal"

# Row 46 (was FraxShare -> now Aave)
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'93.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.01%  "
$ws.Range("E46").Style = "Normal"

# Row 47 (was Aave -> now FraxShare)
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.43%  "
$ws.Range("E47").Style = "Normal"

